$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing C-column values (processing time) for several rows
$ws.Range("C2").Value = 0.016757
$ws.Range("C3").Value = 0.030951
$ws.Range("C4").Value = 0.000501
$ws.Range("C5").Value = 0.037984
$ws.Range("C10").Value = 0.010388
$ws.Range("C11").Value = 0.001081
$ws.Range("C12").Value = 0.000726
$ws.Range("C16").Value = 0.004504
$ws.Range("C17").Value = 0.001329
$ws.Range("C18").Value = 0.001072
$ws.Range("C19").Value = 0.00093
$ws.Range("C20").Value = 0.000992
$ws.Range("C21").Value = 0.001013
$ws.Range("C22").Value = 0.001347
$ws.Range("C23").Value = 0.001031
$ws.Range("C24").Value = 0.000997
$ws.Range("C25").Value = 0.000999
$ws.Range("C26").Value = 0.002
$ws.Range("C27").Value = 0.001293
$ws.Range("C28").Value = 0.000478
$ws.Range("C31").Value = 0
$ws.Range("C38").Value = 0.013909
$ws.Range("C42").Value = 0
$ws.Range("C48").Value = 0.014129
$ws.Range("C49").Value = 0.015671
$ws.Range("C50").Value = 0.015681
$ws.Range("C51").Value = 0.015636
$ws.Range("C52").Value = 0.015617
$ws.Range("C53").Value = 0.015613
$ws.Range("C54").Value = 0.015632
$ws.Range("C55").Value = 0.015627
$ws.Range("C56").Value = 0.01563
$ws.Range("C57").Value = 0.015564
$ws.Range("C58").Value = 0.015676
$ws.Range("C59").Value = 0
$ws.Range("C60").Value = 0.015628
$ws.Range("C61").Value = 0.015625
$ws.Range("C62").Value = 0.015623
$ws.Range("C63").Value = 0.025469
$ws.Range("C64").Value = 0.012012

# Append new rows 65-84 for white_fred_1.jpg .. white_fred_20.jpg
$newRows = @(
    @("white_fred_1.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015687, 0),
    @("white_fred_2.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015633, 0),
    @("white_fred_3.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.025223, 0),
    @("white_fred_4.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.013146, 0),
    @("white_fred_5.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.013916, 0),
    @("white_fred_6.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.020742, 0),
    @("white_fred_7.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.014019, 0),
    @("white_fred_8.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.013901, 0),
    @("white_fred_9.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.013836, 0),
    @("white_fred_10.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.021632, 0),
    @("white_fred_11.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.013031, 0),
    @("white_fred_12.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015703, 0),
    @("white_fred_13.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @("white_fred_14.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.029823, 0),
    @("white_fred_15.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.010093, 0),
    @("white_fred_16.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015694, 0),
    @("white_fred_17.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015623, 0),
    @("white_fred_18.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @("white_fred_19.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015638, 0),
    @("white_fred_20.jpg", "1011100100000000110001000000000000000000000000000001000000000000", 0.015624, 0)
)

$startRow = 65
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
